$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$I = @(7,7,7,5,5,8,6,8,7,9,7,9,4,7,6,7,7,5,6,10,8,9,5,5,7,6,6,7,9,5,8,8,8,7,7,7,11,7,6,6,7,6,6,7,8,7,8,5,3,5,8,6,6,4,4)
$J = @(7,7,7,7,6,8,6,9,7,10,7,9,4,7,7,7,7,6,7,10,8,9,6,6,7,6,8,7,9,6,8,8,8,7,7,7,11,7,7,7,7,7,7,9,8,7,9,6,4,7,9,6,6,4,4)

for ($r = 2; $r -le 56; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $I[$idx]
    $ws.Cells.Item($r, 10).Value = $J[$idx]
}
